$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.816.92"
$ws.Range("D2").Style = "Normal"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.117.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.36%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.77%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5329"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4424"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09016"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "47.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.179"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.00"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.114.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.06%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.773"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.831"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.16%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.18%  "

$ws.Range("E17").Value = "  -0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001135"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.19%  "

$ws.Range("E20").Value = "  +2.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.336"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.46%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.874.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.89%  "

$ws.Range("E24").Value = "  +6.58%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.362.58"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.17%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.278"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.21%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.596"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +10.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "163.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.39%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.192"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1088"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.53%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.228"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.83%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.022"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.43%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.557"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +18.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02606"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.21%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.551"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.24%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "12.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06774"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.93%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.589"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2308"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6857"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.35%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.247"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.85%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6462"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.82%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.86%  "

$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.238"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.00%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.662"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("E49").Value = "  +3.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.03"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.81%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.193"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.71%  "

